$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at M (shifts category/flag/blank header cells from M/N/O to N/O/P,
# and shifts the mrp/category-value/flag-value data cells from L/M/N to M/N/O).
$ws.Columns("M:M").Insert()

# New "price" header (first brand-new shared string -> becomes index 21)
$ws.Range("M1").Value = "price"

# Row 2: the mrp value (still sitting in L2 after the column insert) needs to move to M2,
# and a new price value of 1 goes into L2.
$ws.Range("M2").Value = $ws.Range("L2").Value2
$ws.Range("L2").Value = 1

# Row 2: the flag value's text changes (second brand-new shared string -> becomes index 22)
$ws.Range("O2").Value = "Available or N/A (only Enter 1 from these)"

# Match the saved selection state
$ws.Range("A1:P2").Select()
$excel.ActiveWindow.RangeSelection.Item(1).Activate() | Out-Null
$ws.Range("P2").Activate()
